$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.021.61'
$ws.Range("E2").Value = '  +1.35%  '

$ws.Range("D3").Value = '1.853.25'
$ws.Range("E3").Value = '  +2.75%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = "'236.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.60%  '

$ws.Range("D6").Value = "'0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.65%  '

$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("D8").Value = "'41.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.10%  '

$ws.Range("D9").Value = "'0.328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.00%  '

$ws.Range("E10").Value = '  +2.14%  '

$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("D12").Value = '2.123.62'
$ws.Range("E12").Value = '  +2.95%  '

$ws.Range("E13").Value = '  +2.80%  '

$ws.Range("D14").Value = '1.847.95'
$ws.Range("E14").Value = '  +2.41%  '

$ws.Range("E15").Value = '  +2.39%  '

$ws.Range("D16").Value = "'4.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.83%  '

$ws.Range("D17").Value = '35.010.06'
$ws.Range("E17").Value = '  +1.53%  '

$ws.Range("D18").Value = "'70.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("E19").Value = '  +1.91%  '

$ws.Range("D20").Value = "'240.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.32%  '

$ws.Range("E21").Value = '  +2.69%  '

$ws.Range("D22").Value = "'4.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.78%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("E24").Value = '  +1.74%  '

$ws.Range("D25").Value = "'170.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.50%  '

$ws.Range("D26").Value = "'1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +28.87%  '

$ws.Range("D27").Value = "'7.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.91%  '

$ws.Range("D28").Value = "'17.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.57%  '

$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("E31").Value = '  +2.13%  '

$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").Value = "'4.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.05%  '

$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = "'1.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +23.54%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = "'2.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.78%  '

$ws.Range("D36").Value = "'1.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.09%  '

$ws.Range("D37").Value = "'0.779"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.62%  '

$ws.Range("E38").Value = '  +12.28%  '

$ws.Range("D39").Value = "'91.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.15%  '

$ws.Range("E40").Value = '  +6.73%  '

$ws.Range("D41").Value = '1.349.61'
$ws.Range("E41").Value = '  +1.98%  '

$ws.Range("D42").Value = "'14.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.42%  '

$ws.Range("D43").Value = "'2.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.86%  '

$ws.Range("D44").Value = "'12.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +54.67%  '

$ws.Range("D45").Value = "'2.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.34%  '

$ws.Range("D46").Value = "'0.0554"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.12%  '

$ws.Range("D47").Value = "'2.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("D48").Value = "'6.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.77%  '

$ws.Range("D49").Value = '2.035.96'
$ws.Range("E49").Value = '  +2.40%  '

$ws.Range("D50").Value = "'0.0681"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.41%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = "'3.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.22%  '
